# Adds a new column J ("link" / movie-id slug) across the 2016 movie gross
# list (rows 2-101), mirroring a spreadsheet IMPORTXML/DUMMYFUNCTION lookup
# column that was added upstream. J2 additionally gets a distinct font/fill
# (Inconsolata 11pt on a white fill) matching the header-like styling used
# for the top "link" cell, while J3:J101 stay in the default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - unique lookup formula (top-ranked movie), styled distinctly
$ws.Range("J2").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("importxml(""https://www.boxofficemojo.com/yearly/chart/?view=releasedate&view2=domestic&yr=2016&sort=gross&order=DESC&p=.htm"", ""//*[@id=''body'']/table[3]/tr/td[1]/table[1]/tr/td[2]/table[1]/tr/td/table[1]/tr/td/table[1]/tr/td[2]/b/a/@href"")"),"/movies/?id=starwars2016.htm")'
$ws.Range("J2").Font.Name = "Inconsolata"
$ws.Range("J2").Font.Size = 11
$ws.Range("J2").Interior.Color = 16777215

# Rows 3-101 - plain cached lookup values (default style)
$ws.Range("J3").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pixar2015.htm")'
$ws.Range("J4").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=marvel2016.htm")'
$ws.Range("J5").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=illumination2015.htm")'
$ws.Range("J6").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=junglebook2015.htm")'
$ws.Range("J7").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=deadpool2016.htm")'
$ws.Range("J8").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=disney2016.htm")'
$ws.Range("J9").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=superman2015.htm")'
$ws.Range("J10").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=dc2016.htm")'
$ws.Range("J11").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=illumination2016.htm")'
$ws.Range("J12").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=disney1116.htm")'
$ws.Range("J13").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=fantasticbeasts.htm")'
$ws.Range("J14").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=marvel716.htm")'
$ws.Range("J15").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hiddenfigures.htm")'
$ws.Range("J16").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=bourne5.htm")'
$ws.Range("J17").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=startrek2016.htm")'
$ws.Range("J18").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=x-men2016.htm")'
$ws.Range("J19").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=trolls.htm")'
$ws.Range("J20").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=lalaland.htm")'
$ws.Range("J21").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=kungfupanda3.htm")'
$ws.Range("J22").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=ghostbusters2016.htm")'
$ws.Range("J23").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=centralintelligence.htm")'
$ws.Range("J24").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=tarzan2016.htm")'
$ws.Range("J25").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=sully.htm")'
$ws.Range("J26").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=untitledlucasmoore.htm")'
$ws.Range("J27").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=angrybirds.htm")'
$ws.Range("J28").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=id42.htm")'
$ws.Range("J29").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=conjuring2.htm")'
$ws.Range("J30").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=arrival2016.htm")'
$ws.Range("J31").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=passengers2016.htm")'
$ws.Range("J32").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=sausageparty.htm")'
$ws.Range("J33").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=themagnificentseven.htm")'
$ws.Range("J34").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=ridealong2.htm")'
$ws.Range("J35").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=dontbreathe.htm")'
$ws.Range("J36").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=peregrine.htm")'
$ws.Range("J37").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=theaccountant.htm")'
$ws.Range("J38").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=tmnt2016.htm")'
$ws.Range("J39").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=purge3.htm")'
$ws.Range("J40").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=alice2.htm")'
$ws.Range("J41").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=petesdragon2016.htm")'
$ws.Range("J42").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=thegirlonthetrain2016.htm")'
$ws.Range("J43").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=booamadeahalloween.htm")'
$ws.Range("J44").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=storks.htm")'
$ws.Range("J45").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=badrobot2016.htm")'
$ws.Range("J46").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=newline0116.htm")'
$ws.Range("J47").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hacksawridge.htm")'
$ws.Range("J48").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=allegiant.htm")'
$ws.Range("J49").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=nowyouseeme2.htm")'
$ws.Range("J50").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=iceage5.htm")'
$ws.Range("J51").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=michelledarnell.htm")'
$ws.Range("J52").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=londonhasfallen.htm")'
$ws.Range("J53").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=miraclesfromheaven.htm")'
$ws.Range("J54").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=deepwaterhorizon.htm")'
$ws.Range("J55").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=whyhim.htm")'
$ws.Range("J56").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mybigfatgreekwedding2.htm")'
$ws.Range("J57").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=jackreacher2.htm")'
$ws.Range("J58").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=fences.htm")'
$ws.Range("J59").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mebeforeyou.htm")'
$ws.Range("J60").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=bfg.htm")'
$ws.Range("J61").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=universalcomedy2016.htm")'
$ws.Range("J62").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=theshallows.htm")'
$ws.Range("J63").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=officechristmasparty.htm")'
$ws.Range("J64").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=assassinscreed.htm")'
$ws.Range("J65").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=barbershop3.htm")'
$ws.Range("J66").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=13hoursthesecretsoldiersofbenghazi.htm")'
$ws.Range("J67").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=lion.htm")'
$ws.Range("J68").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=huntsman.htm")'
$ws.Range("J69").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=kuboandthetwostrings.htm")'
$ws.Range("J70").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=manchesterbythesea.htm")'
$ws.Range("J71").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=warcraft.htm")'
$ws.Range("J72").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=howtobesingle.htm")'
$ws.Range("J73").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mikeanddave.htm")'
$ws.Range("J74").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=armsandthedudes.htm")'
$ws.Range("J75").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=ameyerschristmas.htm")'
$ws.Range("J76").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=moneymonster.htm")'
$ws.Range("J77").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pittcotillard.htm")'
$ws.Range("J78").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=nerve.htm")'
$ws.Range("J79").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=risen.htm")'
$ws.Range("J80").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=niceguys.htm")'
$ws.Range("J81").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=boy2016.htm")'
$ws.Range("J82").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=dirtygrandpa.htm")'
$ws.Range("J83").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=blumhouse3.htm")'
$ws.Range("J84").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=5thwave.htm")'
$ws.Range("J85").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=inferno2015.htm")'
$ws.Range("J86").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mothersday.htm")'
$ws.Range("J87").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=patriotsday.htm")'
$ws.Range("J88").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=godsofegypt.htm")'
$ws.Range("J89").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=collateralbeauty.htm")'
$ws.Range("J90").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hailcaesar.htm")'
$ws.Range("J91").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=whentheboughbreaks.htm")'
$ws.Range("J92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=zoolander2.htm")'
$ws.Range("J93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=moonlight2016.htm")'
$ws.Range("J94").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=finesthours.htm")'
$ws.Range("J95").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=florencefosterjenkins.htm")'
$ws.Range("J96").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hellorhighwater.htm")'
$ws.Range("J97").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=forest.htm")'
$ws.Range("J98").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=benhur2016.htm")'
$ws.Range("J99").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=thewitch.htm")'
$ws.Range("J100").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=bridgetjonessbaby.htm")'
$ws.Range("J101").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=kevinhart2016.htm")'
